$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Invalid column (G): rows 3 and 4 become 1
$ws.Range("G3:G4").Value = 1

# Absent column (H): rows 3 through 18 become 1
$ws.Range("H3:H18").Value = 1
